$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 361, shifting rows 361:438 down to 362:439
$ws.Rows.Item(361).Insert()

# Populate the new row 361 with its data
$ws.Cells.Item(361, 1).Value = 10
$ws.Cells.Item(361, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(361, 3).Value = "La Araucanía"
$ws.Cells.Item(361, 4).Value = 45211
$ws.Cells.Item(361, 5).Value = 9
$ws.Cells.Item(361, 6).Value = 100112052
$ws.Cells.Item(361, 7).Value = "Albahaca"
$ws.Cells.Item(361, 8).Value = "Sin especificar"
$ws.Cells.Item(361, 9).Value = "Primera"
$ws.Cells.Item(361, 10).Value = 120
$ws.Cells.Item(361, 11).Value = 6000
$ws.Cells.Item(361, 12).Value = 6000
$ws.Cells.Item(361, 13).Value = 6000
$ws.Cells.Item(361, 14).Value = "$/paquete"
$ws.Cells.Item(361, 15).Value = "Región Metropolitana"
$ws.Cells.Item(361, 16).Value = 6000
$ws.Cells.Item(361, 17).Value = 1
$ws.Cells.Item(361, 18).Value = "Hortaliza"
